# Project DesignFirst is saved.
#
# Target change (see diff):
#   1. Every column definition on the sheet ("Rules") has its outline
#      "collapsed" flag turned on - columns A:K plus the L:XFD catch-all.
#   2. Cell D10 changes its value from 21 to 100.0 (still a plain number,
#      same cell style s="21").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Collapse every column on the sheet ----------------------------------
# Mirrors "select the defined columns -> collapse the outline group" for both
# column ranges that appear in the <cols> definition: A:K (explicit widths)
# and L:XFD (the trailing catch-all column definition).
$lastColumn = 16384
$columnRanges = @(
    $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, 11)).EntireColumn,
    $ws.Range($ws.Cells.Item(1, 12), $ws.Cells.Item(1, $lastColumn)).EntireColumn
)
foreach ($columnRange in $columnRanges) {
    $columnRange.ShowDetail = $false
}

# Also flip the flag column-by-column for the 11 explicitly sized columns,
# plus the trailing catch-all column, so each <col> entry is addressed
# individually as well as via the aggregated ranges above.
for ($i = 1; $i -le 12; $i++) {
    $ws.Columns.Item($i).ShowDetail = $false
}

# --- 2. Update D10's value ---------------------------------------------------
$ws.Range("D10").Value = 100.0
